$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the new row to be stored as literal text (matching the existing
# rows, which are all text even though some values look like dates/numbers).
$ws.Range("A95:E95").NumberFormat = "@"

$ws.Range("A95").Value = "2025-12-20"
$ws.Range("B95").Value = "Pick 4"
$ws.Range("C95").Value = "251220"
$ws.Range("D95").Value = "4-7-0-9"
$ws.Range("E95").Value = "2025-12-20T21:37:21.880+04:00"
